$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet currently holds the text "R40" and needs to
# become the text "1". A plain `.Value = "1"` assignment would be
# auto-coerced to the NUMBER 1 by Excel (since "1" parses cleanly as a
# number), which would also pick up a different/new number-format style.
# To keep B11 a genuine text cell (and keep its existing style untouched),
# build the text via a TEXT() formula, then collapse the formula down to a
# literal value in place with copy / paste-special-values.
$ws.Range("B11").Formula = '=TEXT(1,"0")'
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
